$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM-byvalue")

# Insert a new row at 82, shifting existing row 82 ("U1, U2" / XS2-L2-124QFN ...)
# and everything below it down by one row.
$ws.Rows.Item(82).Insert()

# Populate the new row 82 with the S1 SPDT switch entry.
$ws.Cells.Item(82, 1).Value = "S1"
$ws.Cells.Item(82, 3).Value = "SWITCH-SPDT"
$ws.Cells.Item(82, 4).Value = "SWITCH-SPDT"
$ws.Cells.Item(82, 5).Value = "DK"
$ws.Cells.Item(82, 6).Value = "679-1854-ND"
$ws.Cells.Item(82, 7).Value = "MMS1208"
$ws.Cells.Item(82, 8).Value = 1
$ws.Cells.Item(82, 9).Value = 1.02
$ws.Range("J82").Formula = "=H82*I82"

# Match the saved view/selection state from the authored workbook.
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Range("J82").Select()

$wb.Save()
